# Commit being replicated: the "Files" cell for the second sample row
# (I2 on Sheet1) is updated from the placeholder
#   "test1.fastq.gz, test2.fastq.gz"
# to the real paired-end filenames
#   "test1_R1.fastq.gz, test1_R2.fastq.gz"
# and the sheet's active/selected cell moves from A2 to I2 (the cell the
# author just edited) to match the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Files value for the second row.
$ws.Range("I2").Value = "test1_R1.fastq.gz, test1_R2.fastq.gz"

# Leave this sheet active, with I2 as the selected/active cell (matches
# the <selection pane="bottomRight" activeCell="I2" .../> in the saved
# file) without disturbing the existing frozen-pane setup.
$ws.Activate()
$ws.Range("I2").Select()
